$wb = $excel.ActiveWorkbook
$srcWs = $wb.Worksheets.Item("TwoxTwowAuxDem")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "two_by_two_AuxinInput"

# Column A labels (filled first, top to bottom, to mirror original authoring order)
$ws.Range("A2").Value = 'X'
$ws.Range("A3").Value = 'Y'
$ws.Range("A4").Value = 'W'
$ws.Range("A5").Value = 'TL'
$ws.Range("A6").Value = 'TK'
$ws.Range("A7").Value = 'PX'
$ws.Range("A8").Value = 'PY'
$ws.Range("A9").Value = 'PW'
$ws.Range("A10").Value = 'PL'
$ws.Range("A11").Value = 'PK'
$ws.Range("A12").Value = 'PKS'
$ws.Range("A13").Value = 'PLS'
$ws.Range("A14").Value = 'TAU'
$ws.Range("A15").Value = 'SXX'
$ws.Range("A16").Value = 'SYY'
$ws.Range("A17").Value = 'SWW'
$ws.Range("A18").Value = 'SLSTL'
$ws.Range("A19").Value = 'SKSTK'
$ws.Range("A20").Value = 'DLSX'
$ws.Range("A21").Value = 'DKSX'
$ws.Range("A22").Value = 'DLSY'
$ws.Range("A23").Value = 'DKSY'
$ws.Range("A24").Value = 'DXW'
$ws.Range("A25").Value = 'DYW'
$ws.Range("A26").Value = 'DLW'
$ws.Range("A27").Value = 'DLTL'
$ws.Range("A28").Value = 'DKTK'
$ws.Range("A29").Value = 'CONS'
$ws.Range("A30").Value = 'CWCONS'

# Numeric data values B2:G30
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1.006387575650056
$ws.Range("D2").Value = 1.0126852580817098
$ws.Range("E2").Value = 1.0193022237638398
$ws.Range("F2").Value = 1.026928314669364
$ws.Range("G2").Value = 1.0377545680283367
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1.0095966476462173
$ws.Range("D3").Value = 1.0190881035617243
$ws.Range("E3").Value = 1.0290926053375762
$ws.Range("F3").Value = 1.0406631890486693
$ws.Range("G3").Value = 1.0571630630829592
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.0008350269671129
$ws.Range("D4").Value = 1.0014491088570601
$ws.Range("E4").Value = 1.001862188438271
$ws.Range("F4").Value = 1.0020312308441757
$ws.Range("G4").Value = 1.0016777194890492
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1.0160455226090923
$ws.Range("D5").Value = 1.0320154991633979
$ws.Range("E5").Value = 1.0489563810795042
$ws.Range("F5").Value = 1.0686864948681873
$ws.Range("G5").Value = 1.0970757978651695
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.99062365945204522
$ws.Range("D7").Value = 0.98023361021433086
$ws.Range("E7").Value = 0.96816873805100456
$ws.Range("F7").Value = 0.95318925131473309
$ws.Range("G7").Value = 0.93155409308673487
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.98747489439651848
$ws.Range("D8").Value = 0.97407488427239663
$ws.Range("E8").Value = 0.95895796214597595
$ws.Range("F8").Value = 0.94060887491218026
$ws.Range("G8").Value = 0.91445165767240888
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.99916566972116849
$ws.Range("D9").Value = 0.99732896542206706
$ws.Range("E9").Value = 0.9943694698123019
$ws.Range("F9").Value = 0.99001673768539966
$ws.Range("G9").Value = 0.98375023335117484
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 1.0237432725840783
$ws.Range("D10").Value = 1.0469473301163175
$ws.Range("E10").Value = 1.0711351761369892
$ws.Range("F10").Value = 1.0991328545245336
$ws.Range("G10").Value = 1.1409855970537344
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 0.95689374406213035
$ws.Range("D11").Value = 0.91107593331151671
$ws.Range("E11").Value = 0.85945629637063958
$ws.Range("F11").Value = 0.79586881245096153
$ws.Range("G11").Value = 0.69919000272678622
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.99695134301753163
$ws.Range("D12").Value = 0.9926681265402667
$ws.Range("E12").Value = 0.98685654767402087
$ws.Range("F12").Value = 0.97885703141359282
$ws.Range("G12").Value = 0.96672451546625382
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0.98120735816782156
$ws.Range("D13").Value = 0.96187327355543928
$ws.Range("E13").Value = 0.94079845975902054
$ws.Range("F13").Value = 0.91594404543711139
$ws.Range("G13").Value = 0.88118297509381738
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1.0009381670421438
$ws.Range("D14").Value = 1.024890126404802
$ws.Range("E14").Value = 1.0796578594751951
$ws.Range("F14").Value = 1.1897677711350092
$ws.Range("G14").Value = 1.464804238667671
$ws.Range("B15").Value = 120
$ws.Range("C15").Value = 120
$ws.Range("D15").Value = 120
$ws.Range("E15").Value = 120
$ws.Range("F15").Value = 120
$ws.Range("G15").Value = 120
$ws.Range("B16").Value = 120
$ws.Range("C16").Value = 120
$ws.Range("D16").Value = 120
$ws.Range("E16").Value = 120
$ws.Range("F16").Value = 120
$ws.Range("G16").Value = 120
$ws.Range("B17").Value = 340
$ws.Range("C17").Value = 340
$ws.Range("D17").Value = 340
$ws.Range("E17").Value = 340
$ws.Range("F17").Value = 340
$ws.Range("G17").Value = 340
$ws.Range("B18").Value = 120
$ws.Range("C18").Value = 120
$ws.Range("D18").Value = 120
$ws.Range("E18").Value = 120
$ws.Range("F18").Value = 119.99999999999999
$ws.Range("G18").Value = 120
$ws.Range("B19").Value = 120
$ws.Range("C19").Value = 120
$ws.Range("D19").Value = 120
$ws.Range("E19").Value = 120
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 120
$ws.Range("B20").Value = 48
$ws.Range("C20").Value = 48.460639087018805
$ws.Range("D20").Value = 48.916228970963118
$ws.Range("E20").Value = 49.396445056204186
$ws.Range("F20").Value = 49.951833074336434
$ws.Range("G20").Value = 50.74382702798205
$ws.Range("B21").Value = 72
$ws.Range("C21").Value = 71.543013588470316
$ws.Range("D21").Value = 71.09810222416634
$ws.Range("E21").Value = 70.63655736385509
$ws.Range("F21").Value = 70.11200194940696
$ws.Range("G21").Value = 69.380566675601443
$ws.Range("B22").Value = 72
$ws.Range("C22").Value = 72.459905446804385
$ws.Range("D22").Value = 72.913338581883423
$ws.Range("E22").Value = 73.389760110996932
$ws.Range("F22").Value = 73.938838656194577
$ws.Range("G22").Value = 74.718328898040255
$ws.Range("B23").Value = 48
$ws.Range("C23").Value = 47.543739484384446
$ws.Range("D23").Value = 47.100932522162822
$ws.Range("E23").Value = 46.643032659100584
$ws.Range("F23").Value = 46.124433443138592
$ws.Range("G23").Value = 45.404537555464387
$ws.Range("B24").Value = 120
$ws.Range("C24").Value = 120.66574992281421
$ws.Range("D24").Value = 121.34638684585452
$ws.Range("E24").Value = 122.08891428703427
$ws.Range("F24").Value = 122.98159375382461
$ws.Range("G24").Value = 124.32197077012238
$ws.Range("B25").Value = 120
$ws.Range("C25").Value = 121.05051727124193
$ws.Range("D25").Value = 122.11361650416337
$ws.Range("E25").Value = 123.26157635813239
$ws.Range("F25").Value = 124.62643762174304
$ws.Range("G25").Value = 126.64708928003867
$ws.Range("B26").Value = 100
$ws.Range("C26").Value = 98.313353437742862
$ws.Range("D26").Value = 96.658381566822541
$ws.Range("E26").Value = 94.927588833650503
$ws.Range("F26").Value = 92.9425627130626
$ws.Range("G26").Value = 90.141188584628551
$ws.Range("B27").Value = 100
$ws.Range("C27").Value = 100
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 100
$ws.Range("G27").Value = 100
$ws.Range("B28").Value = 100
$ws.Range("C28").Value = 100
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 100
$ws.Range("G28").Value = 100
$ws.Range("B29").Value = 340
$ws.Range("C29").Value = 340
$ws.Range("D29").Value = 339.58322924414983
$ws.Range("E29").Value = 338.71519886840156
$ws.Range("F29").Value = 337.28941467454194
$ws.Range("G29").Value = 335.03623469860844
$ws.Range("B30").Value = 340
$ws.Range("C30").Value = 340.28390916881636
$ws.Range("D30").Value = 340.49269701139843
$ws.Range("E30").Value = 340.63314406901054
$ws.Range("F30").Value = 340.69061848670492
$ws.Range("G30").Value = 340.57042462627646

# Header row (B1 benchmark already existed; C1:G1 scenario labels added last)
$ws.Range("B1").Value = 'benchmark'
$ws.Range("C1").Value = 'L.15,K.25'
$ws.Range("D1").Value = 'L.1,K.3'
$ws.Range("E1").Value = 'L.05,K.35'
$ws.Range("F1").Value = 'L.0,K.4'
$ws.Range("G1").Value = 'L-.05,K.45'

# Apply quote-prefix text style (matching source sheet) to label columns
$srcCell = $srcWs.Range("A2")
$srcCell.Copy()
$ws.Range("A2:A31").PasteSpecial(-4122)
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "Sheet two_by_two_AuxinInput created"
